$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 96, shifting existing rows 96-108 down to 97-109.
$ws.Rows.Item(96).Insert()

# Populate the newly inserted row 96 with the new weekly record.
$ws.Range("A96").Value = 10
$ws.Range("B96").Value = "Vega Modelo de Temuco"
$ws.Range("C96").Value = "La Araucanía"
$ws.Range("D96").Value = 45131
$ws.Range("E96").Value = 9
$ws.Range("F96").Value = "Fruta"
$ws.Range("G96").Value = 100108
$ws.Range("H96").Value = "Tropicales y subtropicales"
$ws.Range("I96").Value = 100108007
$ws.Range("J96").Value = "Coco"
$ws.Range("K96").Value = "Sin especificar"
$ws.Range("L96").Value = "Primera"
$ws.Range("M96").Value = 50
$ws.Range("N96").Value = 36000
$ws.Range("O96").Value = 36000
$ws.Range("P96").Value = 36000
$ws.Range("Q96").Value = "$/malla 20 unidades"
$ws.Range("R96").Value = "Perú"
$ws.Range("S96").Value = 1800
$ws.Range("T96").Value = 20
